# Bugfixed evaluation and simulated rt_data for components
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (2-8) with corrected forecast error values
$ws.Range("B2").Value = 0.199694501984579
$ws.Range("C2").Value = 1.011608396774954
$ws.Range("D2").Value = 2.462918637783109
$ws.Range("E2").Value = 1.569368866067856
$ws.Range("F2").Value = 1.615372555534037

$ws.Range("B3").Value = 0.181859565044813
$ws.Range("C3").Value = 1.06583629132514
$ws.Range("D3").Value = 2.354939032676854
$ws.Range("E3").Value = 1.53458106096643
$ws.Range("F3").Value = 1.585987066640765

$ws.Range("B4").Value = 0.392354498814844
$ws.Range("C4").Value = 1.092524195814124
$ws.Range("D4").Value = 3.401579278512719
$ws.Range("E4").Value = 1.84433708375468
$ws.Range("F4").Value = 1.882253164048335

$ws.Range("B5").Value = 0.4189316725569809
$ws.Range("C5").Value = 1.35908918265491
$ws.Range("D5").Value = 4.053754439716902
$ws.Range("E5").Value = 2.013393761715999
$ws.Range("F5").Value = 2.065448077970023

$ws.Range("B6").Value = 0.5131610716377353
$ws.Range("C6").Value = 1.045711216584711
$ws.Range("D6").Value = 2.987355008846897
$ws.Range("E6").Value = 1.728396658422741
$ws.Range("F6").Value = 1.739738397768312

$ws.Range("B7").Value = 0.004917215416235725
$ws.Range("C7").Value = 0.7811863174655274
$ws.Range("D7").Value = 1.08329812457537
$ws.Range("E7").Value = 1.040816085855407
$ws.Range("F7").Value = 1.103939848344968
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = -0.3858263734312975
$ws.Range("C8").Value = 0.875567608838748
$ws.Range("D8").Value = 1.133368546339351
$ws.Range("E8").Value = 1.064597833146091
$ws.Range("F8").Value = 1.086925879296765
$ws.Range("G8").Value = 6

# Row 9: update values and add F9 (previously missing)
$ws.Range("B9").Value = -0.8536607533201513
$ws.Range("C9").Value = 0.8536607533201513
$ws.Range("D9").Value = 0.9321950237660799
$ws.Range("E9").Value = 0.9655024721698438
$ws.Range("F9").Value = 0.5524377910773552
$ws.Range("G9").Value = 3

# Row 10: new row "Q8" (copy formatting from A9, which holds the same style)
$ws.Range("A10").Value = "Q8"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = 0.265633061926664
$ws.Range("C10").Value = 0.265633061926664
$ws.Range("D10").Value = 0.07056092358853493
$ws.Range("E10").Value = 0.265633061926664
$ws.Range("G10").Value = 1
